# Step Size Calculator - add Conversion / Equivalent Steps section
# and change the limit-switch check frequency (B2) + related values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Column A / B / C block (left table)
# ---------------------------------------------------------------

# Row2: Length - change check frequency value (B2) 0.15 -> 0.025
$ws.Range("B2").Value = 0.025

# Row5 stays the same (Step size 0.5)

# Row7 becomes "Conversion" (was "Steps (rounded)"), formula moves here
$ws.Range("A7").Value = "Conversion"
$ws.Range("B7").Formula = "=360/(2*PI()*B3/100*B4*B5)"
$ws.Range("B7").Style = "Normal"

# Row8 (new) becomes "Steps (rounded)" with formula referencing B7*B2
$ws.Range("A8").Value = "Steps (rounded)"
$ws.Range("B8").Formula = "=ROUND(B7*B2,0)"
$ws.Range("B8").Style = "Output"

# ---------------------------------------------------------------
# Column H header - Equivalent Steps (added right after "Conversion")
# ---------------------------------------------------------------

$ws.Range("H1").Value = "Equivalent Steps"

# ---------------------------------------------------------------
# Column E / F / G block (right table) - Key Measurements / Conversion
# ---------------------------------------------------------------

# Row2: Y Body Length 20 [cm] -> 0.2 [m]
$ws.Range("F2").Value = 0.2
$ws.Range("G2").Value = "[m]"

# Row3: X Body Length 18 [cm] -> 0.18 [m]
$ws.Range("F3").Value = 0.18
$ws.Range("G3").Value = "[m]"

# Row4: Gantry Working Size stays 1.55 [m] (unchanged)

# Row5 (new): X Offset
$ws.Range("E5").Value = "X Offset"
$ws.Range("F5").Value = 0.05
$ws.Range("G5").Value = "[m]"

# Row6 (new): Y Offset
$ws.Range("E6").Value = "Y Offset"
$ws.Range("F6").Value = 0.025
$ws.Range("G6").Value = "[m]"

# Row7 (new): Boundary Offset
$ws.Range("E7").Value = "Boundary Offset"
$ws.Range("F7").Value = 0.3
$ws.Range("G7").Value = "[m]"

# ---------------------------------------------------------------
# Column H - Equivalent Steps formulas
# ---------------------------------------------------------------

$ws.Range("H2").Formula = '=ROUND(F2*$B$7,0)'
$ws.Range("H3").Formula = '=ROUND(F3*$B$7,0)'
$ws.Range("H4").Formula = '=ROUND(F4*$B$7,0)'
$ws.Range("H5").Formula = '=ROUND(F5*$B$7,0)'
$ws.Range("H6").Formula = '=ROUND(F6*$B$7,0)'
$ws.Range("H7").Formula = '=ROUND(F7*$B$7,0)'

# ---------------------------------------------------------------
# Sheet view tidy-up (selection moved when rows were inserted)
# ---------------------------------------------------------------

[void]$ws.Range("F8").Select()
